# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2310
#   *_new  -> *_FV2404
# Also: freeze the header row and wrap the data range in an Excel Table
# (ListObject) named "Table1", matching the regenerated AHB export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row ------------------------------------------------
# Columns A:J used the "_old" suffix, columns L:U used the "_new" suffix
# (column K holds the "diff" header and stays untouched).
$oldSuffixCols = 1..10
$newSuffixCols = 12..21

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Cells.Item(1, $col)
    $value = $cell.Value()
    if ($value -like "*_old") {
        $cell.Value = ($value -replace "_old$", "_FV2310")
    }
}

foreach ($col in $newSuffixCols) {
    $cell = $ws.Cells.Item(1, $col)
    $value = $cell.Value()
    if ($value -like "*_new") {
        $cell.Value = ($value -replace "_new$", "_FV2404")
    }
}

# --- 2. Turn the data range A1:U80 into an Excel Table -----------------------
$tableRange = $ws.Range("A1:U80")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""

# --- 3. Freeze the header row -------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
